$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "gender" -> "gender_" and "address,2_location" -> "addressj,2"
$ws.Range("H1").Value = "gender_"
$ws.Range("D1").Value = "addressj,2"

# Update the active selection on the sheet from H5 to D1
$ws.Range("D1").Select() | Out-Null
